$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Modelo" header in F1, copying the header style from E1 (Tipo)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Modelo"

# Update the refreshed metric values for row 2
$ws.Range("B2").Value = 0.1261900591977314
$ws.Range("C2").Value = 0.9906766187195095
$ws.Range("D2").Value = 0.2768148751524423

# Record which model produced these metrics
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5, n_estimators=50))])"
